$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new columns for Anthony's machines (set in this order so the
# shared-string table gets new entries appended as: desktop, laptop)
$ws.Range("H1").Value = "Anthony's desktop"
$ws.Range("G1").Value = "Anthony's laptop"

# Rename the "Daryl's computer" column header to "Daryl's laptop"
$ws.Range("F1").Value = "Daryl's laptop"

# Update the baseline2 row (row 3) with new measurement values
$ws.Range("C3").Value = 121.152892872923
$ws.Range("G3").Value = 104.580126871513

# Move the active selection to reflect the new working cell
$ws.Range("G11").Select()
